$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 52000
$ws.Range("J3").Value = 52000
$ws.Range("L3").Value = 52000
$ws.Range("N3").Value = -52228
$ws.Range("H9").Value = 630
$ws.Range("J9").Value = 750
$ws.Range("L9").Value = 750
$ws.Range("N9").Value = -1088
$ws.Range("H28").Value = 2028.3
$ws.Range("I28").Value = 1740.7693
$ws.Range("J28").Value = 2562.2856
$ws.Range("K28").Value = 1740.7693
$ws.Range("L28").Value = 2562.2856
$ws.Range("M28").Value = -1255.7693
$ws.Range("N28").Value = -3532.2856
$ws.Range("H53").Value = 304.93332
$ws.Range("I53").Value = 537.75
$ws.Range("J53").Value = 220.27272
$ws.Range("K53").Value = 537.75
$ws.Range("L53").Value = 220.27272
$ws.Range("M53").Value = 99.25
$ws.Range("N53").Value = -1494.27272
$ws.Range("H70").Value = 4745.9165
$ws.Range("J70").Value = 6778.857
$ws.Range("L70").Value = 20336.571
$ws.Range("N70").Value = -20876.571
$ws.Range("H73").Value = 4745.9165
$ws.Range("J73").Value = 6778.857
$ws.Range("L73").Value = 20336.571
$ws.Range("N73").Value = -22208.571
$ws.Range("H92").Value = 31250576
$ws.Range("I92").Value = 35714920
$ws.Range("J92").Value = 168.5
$ws.Range("K92").Value = 35714920
$ws.Range("L92").Value = 168.5
$ws.Range("M92").Value = -35713672
$ws.Range("N92").Value = -2664.5
$ws.Range("H102").Value = 52000
$ws.Range("J102").Value = 52000
$ws.Range("L102").Value = 52000
$ws.Range("N102").Value = -58490
$ws.Range("H111").Value = 2710.611
$ws.Range("I111").Value = 1066.1666
$ws.Range("J111").Value = 5999.5
$ws.Range("K111").Value = 3198.4998
$ws.Range("L111").Value = 17998.5
$ws.Range("M111").Value = -131.4998000000001
$ws.Range("N111").Value = -24132.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15622.468
$ws.Range("I32").Value = 15664.06
$ws.Range("K32").Value = 15664.06
$ws.Range("M32").Value = -15377.06
$ws.Range("H53").Value = 2525000
$ws.Range("I53").Value = 2525000
$ws.Range("K53").Value = 2525000
$ws.Range("M53").Value = -2524318
$ws.Range("H61").Value = 12365
$ws.Range("I61").Value = 13080.4
$ws.Range("J61").Value = 6999.5
$ws.Range("K61").Value = 13080.4
$ws.Range("L61").Value = 6999.5
$ws.Range("M61").Value = -12868.4
$ws.Range("N61").Value = -7423.5
$ws.Range("H74").Value = 1475.9032
$ws.Range("I74").Value = 1171.3478
$ws.Range("K74").Value = 1171.3478
$ws.Range("M74").Value = -297.3478
$ws.Range("H77").Value = 1475.9032
$ws.Range("I77").Value = 1171.3478
$ws.Range("K77").Value = 5856.739
$ws.Range("M77").Value = -1488.739
$ws.Range("H110").Value = 5111791
$ws.Range("I110").Value = 6809387.5
$ws.Range("J110").Value = 19000
$ws.Range("K110").Value = 6809387.5
$ws.Range("L110").Value = 19000
$ws.Range("M110").Value = -6807342.5
$ws.Range("N110").Value = -23090
$ws.Range("H122").Value = 3758.4614
$ws.Range("I122").Value = 2078.9644
$ws.Range("K122").Value = 6236.8932
$ws.Range("M122").Value = -3786.8932
$ws.Range("H136").Value = 12365
$ws.Range("I136").Value = 13080.4
$ws.Range("J136").Value = 6999.5
$ws.Range("K136").Value = 39241.2
$ws.Range("L136").Value = 20998.5
$ws.Range("M136").Value = -36691.2
$ws.Range("N136").Value = -26098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 59404.668
$ws.Range("J76").Value = 59404.668
$ws.Range("L76").Value = 59404.668
$ws.Range("N76").Value = -60034.668
$ws.Range("H79").Value = 59404.668
$ws.Range("J79").Value = 59404.668
$ws.Range("L79").Value = 59404.668
$ws.Range("N79").Value = -61588.668
$ws.Range("H134").Value = 1196.381
$ws.Range("I134").Value = 1196.381
$ws.Range("K134").Value = 3589.143
$ws.Range("M134").Value = -1054.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1457.3334
$ws.Range("I22").Value = 774.5
$ws.Range("K22").Value = 774.5
$ws.Range("M22").Value = -424.5
$ws.Range("H31").Value = 6470.5757
$ws.Range("I31").Value = 3060.75
$ws.Range("J31").Value = 7228.315
$ws.Range("K31").Value = 3060.75
$ws.Range("L31").Value = 7228.315
$ws.Range("M31").Value = -2765.75
$ws.Range("N31").Value = -7818.315
$ws.Range("H34").Value = 6470.5757
$ws.Range("I34").Value = 3060.75
$ws.Range("J34").Value = 7228.315
$ws.Range("K34").Value = 3060.75
$ws.Range("L34").Value = 7228.315
$ws.Range("M34").Value = -2858.75
$ws.Range("N34").Value = -7632.315
$ws.Range("H107").Value = 1299273.1
$ws.Range("I107").Value = 1653391.9
$ws.Range("K107").Value = 1653391.9
$ws.Range("M107").Value = -1651471.9
$ws.Range("H122").Value = 4080.75
$ws.Range("J122").Value = 6435.5
$ws.Range("L122").Value = 19306.5
$ws.Range("N122").Value = -24206.5
$ws.Range("H132").Value = 9707.069
$ws.Range("I132").Value = 10000.179
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 30000.537
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -27470.537
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 208.71428
$ws.Range("I14").Value = 208.71428
$ws.Range("K14").Value = 626.14284
$ws.Range("M14").Value = -453.14284
$ws.Range("H42").Value = 5936.5
$ws.Range("J42").Value = 6570.2856
$ws.Range("L42").Value = 19710.8568
$ws.Range("N42").Value = -20778.8568
$ws.Range("H64").Value = 10307.714
$ws.Range("J64").Value = 12230.8
$ws.Range("L64").Value = 36692.39999999999
$ws.Range("N64").Value = -37232.39999999999
$ws.Range("H67").Value = 10307.714
$ws.Range("J67").Value = 12230.8
$ws.Range("L67").Value = 36692.39999999999
$ws.Range("N67").Value = -38564.39999999999
$ws.Range("H92").Value = 1387
$ws.Range("I92").Value = 1497
$ws.Range("K92").Value = 4491
$ws.Range("M92").Value = -3243
$ws.Range("H129").Value = 1128.3636
$ws.Range("I129").Value = 905.6
$ws.Range("J129").Value = 1225.2174
$ws.Range("K129").Value = 2716.8
$ws.Range("L129").Value = 3675.6522
$ws.Range("M129").Value = 2283.2
$ws.Range("N129").Value = -13675.6522
$ws.Range("H131").Value = 14144258
$ws.Range("I131").Value = 1215
$ws.Range("J131").Value = 15154475
$ws.Range("K131").Value = 3645
$ws.Range("L131").Value = 45463425
$ws.Range("M131").Value = 1395
$ws.Range("N131").Value = -45473505
$ws.Range("H132").Value = 4509.933
$ws.Range("I132").Value = 1162.25
$ws.Range("J132").Value = 5727.273
$ws.Range("K132").Value = 10460.25
$ws.Range("L132").Value = 51545.457
$ws.Range("M132").Value = -7930.25
$ws.Range("N132").Value = -56605.457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 486.4643
$ws.Range("I97").Value = 423.125
$ws.Range("J97").Value = 570.9167
$ws.Range("K97").Value = 423.125
$ws.Range("L97").Value = 570.9167
$ws.Range("M97").Value = 72.875
$ws.Range("N97").Value = -1562.9167
$ws.Range("H113").Value = 3777.8
$ws.Range("I113").Value = 3296.3333
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 3296.3333
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -1126.3333
$ws.Range("N113").Value = -8840
$ws.Range("H122").Value = 463633.88
$ws.Range("I122").Value = 919519.3
$ws.Range("J122").Value = 7748.4165
$ws.Range("K122").Value = 2758557.9
$ws.Range("L122").Value = 23245.2495
$ws.Range("M122").Value = -2756107.9
$ws.Range("N122").Value = -28145.2495
$ws.Range("H126").Value = 4609.7
$ws.Range("I126").Value = 2218.375
$ws.Range("J126").Value = 6203.9165
$ws.Range("K126").Value = 6655.125
$ws.Range("L126").Value = 18611.7495
$ws.Range("M126").Value = -4185.125
$ws.Range("N126").Value = -23551.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5768.154
$ws.Range("J46").Value = 6004.0415
$ws.Range("L46").Value = 6004.0415
$ws.Range("N46").Value = -6380.0415
$ws.Range("H61").Value = 5570.2856
$ws.Range("I61").Value = 4665.6665
$ws.Range("J61").Value = 6248.75
$ws.Range("K61").Value = 4665.6665
$ws.Range("L61").Value = 6248.75
$ws.Range("M61").Value = -4463.6665
$ws.Range("N61").Value = -6652.75
$ws.Range("H113").Value = 5570.2856
$ws.Range("I113").Value = 4665.6665
$ws.Range("J113").Value = 6248.75
$ws.Range("K113").Value = 4665.6665
$ws.Range("L113").Value = 6248.75
$ws.Range("M113").Value = -2495.6665
$ws.Range("N113").Value = -10588.75
$ws.Range("H122").Value = 11365.917
$ws.Range("I122").Value = 5107.1665
$ws.Range("K122").Value = 15321.4995
$ws.Range("M122").Value = -12871.4995
$ws.Range("H134").Value = 81619
$ws.Range("J134").Value = 81619
$ws.Range("L134").Value = 81619
$ws.Range("N134").Value = -91759
$ws.Range("H136").Value = 2353.0227
$ws.Range("I136").Value = 1519.8108
$ws.Range("J136").Value = 6757.143
$ws.Range("K136").Value = 4559.4324
$ws.Range("L136").Value = 20271.429
$ws.Range("M136").Value = -2009.4324
$ws.Range("N136").Value = -25371.429
